# week 2 assignment submitted
# Fill in the "Actual time length to complete" column (C) on the week2
# sheet for the two tasks that were just finished: "DQ2 response 4"
# (row 16) got its actual time recorded, and "Hand-in assignment"
# (row 18) had its actual time corrected/updated. The "Total" formula
# in C19 recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DQ2 response 4 - actual time matches the anticipated 15 minutes
$ws.Range("C16").Value2 = 0.010416666666666666

# Hand-in assignment - actual time took 8 hours instead of the
# previously recorded 6 hours
$ws.Range("C18").Value2 = 0.33333333333333331

# Leave the cursor/selection on the Total row, reflecting where the
# user ended up after entering the data
$ws.Range("C19").Select() | Out-Null
